$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.914.94'
$ws.Range("E2").Value = '  +8.07%  '

$ws.Range("D3").Value = '1.819.76'
$ws.Range("E3").Value = '  +5.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4940'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.77'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2782'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06402'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.71%  '

$ws.Range("D11").Value = '1.808.64'
$ws.Range("E11").Value = '  +4.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.75'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07063'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.13%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.53%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6430'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.672'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.78%  '

$ws.Range("D17").Value = '28.937.96'
$ws.Range("E17").Value = '  +8.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9993'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007316'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9991'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.64%  '

$ws.Range("D22").Value = '2.041.35'

$ws.Range("E23").Value = '  +3.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.833'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.352'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '128.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +21.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.881'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.406'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.133'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08366'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.787'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04920'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.71%  '

$ws.Range("E35").Value = '  +9.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6742'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.27%  '

$ws.Range("E37").Value = '  +4.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.300'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.707'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9482'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.149'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01586'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9993'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4082'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.178'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.16%  '

$ws.Range("E47").Value = '  +5.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05518'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.65'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.090'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.303'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.67%  '
